# Update "想去人数" (interested-count) figures in column F across the
# "展览" and "全部类型" sheets, and the "本地生活" sheet, to match the
# newly scraped data as of commit 456a3b4.

$wb = $excel.ActiveWorkbook

# ---- Sheet: 展览 (Exhibitions) ----
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value  = 2627
$ws1.Range("F7").Value  = 1936
$ws1.Range("F8").Value  = 1773
$ws1.Range("F11").Value = 2432
$ws1.Range("F13").Value = 208
$ws1.Range("F18").Value = 8939
$ws1.Range("F20").Value = 6954
$ws1.Range("F21").Value = 11349
$ws1.Range("F27").Value = 2482
$ws1.Range("F30").Value = 2392
$ws1.Range("F31").Value = 591
$ws1.Range("F33").Value = 4480
$ws1.Range("F34").Value = 757
$ws1.Range("F35").Value = 326
$ws1.Range("F37").Value = 487

# ---- Sheet: 本地生活 (Local life) ----
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F3").Value = 620

# ---- Sheet: 全部类型 (All types) ----
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value  = 620
$ws4.Range("F6").Value  = 2627
$ws4.Range("F10").Value = 1936
$ws4.Range("F12").Value = 1773
$ws4.Range("F15").Value = 2432
$ws4.Range("F18").Value = 208
$ws4.Range("F23").Value = 8939
$ws4.Range("F25").Value = 6954
$ws4.Range("F26").Value = 11349
$ws4.Range("F34").Value = 2482
$ws4.Range("F40").Value = 4480
$ws4.Range("F46").Value = 487
